$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 6.371
$ws.Range("A12").Value = -21.576
$ws.Range("B12").Value = 6.695
$ws.Range("B14").Value = 6.532000000000001
$ws.Range("B22").Value = 6.65
